$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# 1. Append a 4th bullet point to the "Edited corpus" note (cell C9)
$ws.Range("C9").Value = "Edited corpus:`n1. . -> _`n2. Changed pronouns to the referred person`n3. Joined sentences that talked about same thing`n4. Deleted names of actors inside ()"

# 2. Grow row 9 to fit the extra line of text
$ws.Rows.Item(9).RowHeight = 85

# 3. Move the active selection on the Log sheet to C13
$ws.Range("C13").Select()
